$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "CURRENT AS OF:" label moves from D3.2 to D4.1
$ws.Range("R1").Value = "D4.1"

# Move the active selection to A5
$ws.Range("A5").Select()

# Record ordnance expended during period D4.1 (column J)
$ws.Range("J12").Value = 2
$ws.Range("J15").Value = 9
$ws.Range("J18").Value = 3
$ws.Range("J19").Value = 4
